$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# general_single_cell
$ws.Range("C2").Value = 50

# general_column_cells
$ws.Range("C4").Value = 20
$ws.Range("C5").Value = 15
$ws.Range("C6").Value = 20
$ws.Range("C7").Value = 25
$ws.Range("C8").Value = 30
$ws.Range("C9").Value = 35

# general_row_cells
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 6
